$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# 1) "What is Lambda?" heading paragraph -> bold (w:b / w:bCs) on both the
#    paragraph mark run-properties and the run itself.
$pLambda = $d.Paragraphs.Item(2)
if ($pLambda.Range.Text.TrimEnd() -eq "What is Lambda?") {
    $xmlLambda = '<w:p ' + $wns + ' w:rsidR="005B4312" w:rsidRDefault="005B4312"><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>What is Lambda?</w:t></w:r></w:p>'
    [void]$pLambda.Range.InsertXML($xmlLambda)
} else {
    Write-Output "WARN: paragraph 2 text mismatch: [$($pLambda.Range.Text)]"
}

# 2) "Create table in DynampDB." heading paragraph -> bold all three runs
#    (including the misspelled word wrapped in proofErr tags) plus the
#    paragraph mark run-properties.
$pDynamo = $d.Paragraphs.Item(10)
if ($pDynamo.Range.Text.TrimEnd() -eq "Create table in DynampDB.") {
    $xmlDynamo = '<w:p ' + $wns + ' w:rsidR="005A142B" w:rsidRDefault="005A142B" w:rsidP="005B4312"><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Create table in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>DynampDB</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r></w:p>'
    [void]$pDynamo.Range.InsertXML($xmlDynamo)
} else {
    Write-Output "WARN: paragraph 10 text mismatch: [$($pDynamo.Range.Text)]"
}

# 3) "Services search for IAM" -> split into "...IA" + "M" with a _GoBack
#    bookmark inserted between them (the _GoBack bookmark moved here from
#    further down the document - see step 4).
$pIAM = $d.Paragraphs.Item(15)
if ($pIAM.Range.Text.TrimEnd() -eq "Services search for IAM") {
    $xmlIAM = '<w:p ' + $wns + ' w:rsidR="005A142B" w:rsidRDefault="00432E50" w:rsidP="005A142B"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Services search for IA</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>M</w:t></w:r></w:p>'
    [void]$pIAM.Range.InsertXML($xmlIAM)
} else {
    Write-Output "WARN: paragraph 15 text mismatch: [$($pIAM.Range.Text)]"
}

# 4) Remove the old _GoBack bookmark that used to sit after "...based upon
#    requirement." (it effectively moved to step 3 above).
$pPermission = $d.Paragraphs.Item(18)
if ($pPermission.Range.Text.TrimEnd() -eq "To give permission > select DynamoDB and give permissions based upon requirement.") {
    $xmlPermission = '<w:p ' + $wns + ' w:rsidR="00432E50" w:rsidRPr="005A142B" w:rsidRDefault="00432E50" w:rsidP="005A142B"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>To give permission &gt; select DynamoDB and give permissions based upon requirement.</w:t></w:r></w:p>'
    [void]$pPermission.Range.InsertXML($xmlPermission)
} else {
    Write-Output "WARN: paragraph 18 text mismatch: [$($pPermission.Range.Text)]"
}
